# Hortaliza, Vega Modelo de Temuco - Betarraga
# Insert a new weekly record at row 215 (pushing the existing rows 215-252
# down to 216-253), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 215:252 down by one row, creating a blank row 215.
$ws.Rows.Item(215).Insert()

# Fill in the new row 215 with the new weekly observation.
$ws.Cells.Item(215, 1).Value = 10
$ws.Cells.Item(215, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(215, 3).Value = "La Araucanía"
$ws.Cells.Item(215, 4).Value = 44474
$ws.Cells.Item(215, 5).Value = 9
$ws.Cells.Item(215, 6).Value = 100114014
$ws.Cells.Item(215, 7).Value = "Betarraga"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 40
$ws.Cells.Item(215, 11).Value = 8000
$ws.Cells.Item(215, 12).Value = 9000
$ws.Cells.Item(215, 13).Value = 8500
$ws.Cells.Item(215, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(215, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(215, 16).Value = 708
$ws.Cells.Item(215, 17).Value = 12
$ws.Cells.Item(215, 18).Value = "Hortaliza"
